# Workbook / worksheet references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "Förändrad" (changed) date in column C for every data row ---
# All data rows run from row 2 to row 295 and previously held serial date 45184
# (2023-09-15); they must be updated to serial date 45186 (2023-09-17).
$ws.Range("C2:C295").Value = 45186

# --- 2) Add a friendly display-text second argument to the HYPERLINK() formulas
#        in columns S, T, V, W, X, Y for the rows that contain them (rows 2-7). ---
$linkColumns = @(
    @{ Col = "S"; Folder = "artfynd";        Ext = ".xlsx" },
    @{ Col = "T"; Folder = "kartor";         Ext = ".png"  },
    @{ Col = "V"; Folder = "klagomål";       Ext = ".docx" },
    @{ Col = "W"; Folder = "klagomålsmail";  Ext = ".docx" },
    @{ Col = "X"; Folder = "tillsyn";        Ext = ".docx" },
    @{ Col = "Y"; Folder = "tillsynsmail";   Ext = ".docx" }
)

for ($row = 2; $row -le 7; $row++) {
    $beteckning = $ws.Cells.Item($row, 1).Value()

    foreach ($entry in $linkColumns) {
        $col = $entry.Col
        $url = "https://klasma.github.io/Logging_TIDAHOLM/" + $entry.Folder + "/" + $beteckning + $entry.Ext
        $formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        $ws.Range($col + $row).Formula = $formula
    }
}
